# Inserts a new weekly price record for Espinaca (Femacal de La Calera)
# as row 158, pushing all subsequent rows (158-227) down by one
# (new last row becomes 228), matching the "Fruta / hortaliza, semanal"
# update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 158..227 down to 159..228
$ws.Rows.Item(158).Insert()

# Populate the new row 158 with the newly reported weekly data
$ws.Cells.Item(158, 1).Value = 3
$ws.Cells.Item(158, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(158, 3).Value = "Coquimbo"
$ws.Cells.Item(158, 4).Value = 44523
$ws.Cells.Item(158, 5).Value = 5
$ws.Cells.Item(158, 6).Value = 100112012
$ws.Cells.Item(158, 7).Value = "Espinaca"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 160
$ws.Cells.Item(158, 11).Value = 3000
$ws.Cells.Item(158, 12).Value = 3000
$ws.Cells.Item(158, 13).Value = 3000
$ws.Cells.Item(158, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(158, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(158, 16).Value = 1000
$ws.Cells.Item(158, 17).Value = 3
$ws.Cells.Item(158, 18).Value = "Hortaliza"
